$wb = $excel.ActiveWorkbook

# Update "想去人数" (want-to-go count) figures in column F for rows 3-6
# on both the "展览" and "全部类型" sheets, which carry duplicate data.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 436
    $ws.Range("F4").Value = 3184
    $ws.Range("F5").Value = 79
    $ws.Range("F6").Value = 643
}
